# Additional companies sent for questionaire
# Remove the "Parent company" (column B) and "Location County/City" (column E)
# columns from the known locomotive list, shifting the remaining columns left.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column E first (Location County/City) so column B's index isn't
# affected, then delete column B (Parent company).
$ws.Columns("E").Delete()
$ws.Columns("B").Delete()

# Restore the selection that was active in the saved workbook.
$ws.Range("I10").Select()
